# Auto-generated Excel COM-interop edit script
# Applies scheduled-runner price/profit updates to the Sheets workbook.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 505.14285
$ws.Range("I2").Value = 484.6316
$ws.Range("K2").Value = 484.6316
$ws.Range("M2").Value = -371.6316
$ws.Range("H9").Value = 147.40909
$ws.Range("J9").Value = 156.4
$ws.Range("L9").Value = 156.4
$ws.Range("N9").Value = -494.4
$ws.Range("H40").Value = 4999.5
$ws.Range("I40").Value = 4999
$ws.Range("K40").Value = 4999
$ws.Range("M40").Value = -4824
$ws.Range("H100").Value = 14770729
$ws.Range("I100").Value = 17900486
$ws.Range("K100").Value = 17900486
$ws.Range("M100").Value = -17899945
$ws.Range("H107").Value = 6840.24
$ws.Range("I107").Value = 8807.200000000001
$ws.Range("J107").Value = 3889.8
$ws.Range("K107").Value = 8807.200000000001
$ws.Range("L107").Value = 3889.8
$ws.Range("M107").Value = -6887.200000000001
$ws.Range("N107").Value = -7729.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 6894.2173
$ws.Range("I2").Value = 8358.294
$ws.Range("K2").Value = 8358.294
$ws.Range("M2").Value = -8245.294
$ws.Range("H4").Value = 1999
$ws.Range("J4").Value = 3000
$ws.Range("L4").Value = 3000
$ws.Range("N4").Value = -3232
$ws.Range("H34").Value = 300000
$ws.Range("H45").Value = 5704.0713
$ws.Range("I45").Value = 4214.273
$ws.Range("K45").Value = 4214.273
$ws.Range("M45").Value = -3837.273
$ws.Range("H86").Value = 27642
$ws.Range("I86").Value = 27642
$ws.Range("K86").Value = 27642
$ws.Range("M86").Value = -26456
$ws.Range("H89").Value = 27642
$ws.Range("I89").Value = 27642
$ws.Range("K89").Value = 82926
$ws.Range("M89").Value = -76998
$ws.Range("H116").Value = 6894.2173
$ws.Range("I116").Value = 8358.294
$ws.Range("K116").Value = 8358.294
$ws.Range("M116").Value = -6064.294

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 6894.2173
$ws.Range("I3").Value = 8358.294
$ws.Range("K3").Value = 8358.294
$ws.Range("M3").Value = -8244.294
$ws.Range("H80").Value = 798.3570999999999
$ws.Range("J80").Value = 782.8461
$ws.Range("L80").Value = 782.8461
$ws.Range("N80").Value = -2778.8461
$ws.Range("H83").Value = 798.3570999999999
$ws.Range("J83").Value = 782.8461
$ws.Range("L83").Value = 3914.2305
$ws.Range("N83").Value = -13898.2305
$ws.Range("H107").Value = 5192.346
$ws.Range("I107").Value = 5147.8696
$ws.Range("J107").Value = 5533.3335
$ws.Range("K107").Value = 5147.8696
$ws.Range("L107").Value = 5533.3335
$ws.Range("M107").Value = -3227.8696
$ws.Range("N107").Value = -9373.333500000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 125573.8
$ws.Range("J62").Value = 202790.67
$ws.Range("L62").Value = 202790.67
$ws.Range("N62").Value = -204038.67
$ws.Range("H65").Value = 125573.8
$ws.Range("J65").Value = 202790.67
$ws.Range("L65").Value = 1013953.35
$ws.Range("N65").Value = -1020193.35
$ws.Range("H97").Value = 76133
$ws.Range("J97").Value = 76133
$ws.Range("L97").Value = 76133
$ws.Range("N97").Value = -78115
$ws.Range("H99").Value = 226486.27
$ws.Range("I99").Value = 429599.5
$ws.Range("K99").Value = 429599.5
$ws.Range("M99").Value = -428101.5
$ws.Range("H126").Value = 226486.27
$ws.Range("I126").Value = 429599.5
$ws.Range("K126").Value = 1288798.5
$ws.Range("M126").Value = -1286328.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H59").Value = 2875
$ws.Range("I59").Value = 2875
$ws.Range("K59").Value = 8625
$ws.Range("M59").Value = -8085
$ws.Range("H75").Value = 240
$ws.Range("I75").Value = 254
$ws.Range("J75").Value = 212
$ws.Range("K75").Value = 762
$ws.Range("L75").Value = 636
$ws.Range("M75").Value = 236
$ws.Range("N75").Value = -2632
$ws.Range("H78").Value = 240
$ws.Range("I78").Value = 254
$ws.Range("J78").Value = 212
$ws.Range("K78").Value = 2286
$ws.Range("L78").Value = 1908
$ws.Range("M78").Value = 2706
$ws.Range("N78").Value = -11892
$ws.Range("H86").Value = 1097.3334
$ws.Range("J86").Value = 1063.909
$ws.Range("L86").Value = 3191.727
$ws.Range("N86").Value = -5563.727000000001
$ws.Range("H89").Value = 1097.3334
$ws.Range("J89").Value = 1063.909
$ws.Range("L89").Value = 9575.181
$ws.Range("N89").Value = -21431.181
$ws.Range("H132").Value = 57320.89
$ws.Range("I132").Value = 777.8
$ws.Range("J132").Value = 127999.75
$ws.Range("K132").Value = 7000.2
$ws.Range("L132").Value = 1151997.75
$ws.Range("M132").Value = -4470.2
$ws.Range("N132").Value = -1157057.75
$ws.Range("H140").Value = 9846.237999999999
$ws.Range("I140").Value = 10540.579
$ws.Range("K140").Value = 31621.737
$ws.Range("M140").Value = -26441.737

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H33").Value = 19994.5
$ws.Range("J33").Value = 19994.5
$ws.Range("L33").Value = 19994.5
$ws.Range("N33").Value = -20498.5
$ws.Range("H97").Value = 11096.608
$ws.Range("J97").Value = 4000
$ws.Range("L97").Value = 4000
$ws.Range("N97").Value = -4992
$ws.Range("H126").Value = 12297
$ws.Range("I126").Value = 16900.834
$ws.Range("J126").Value = 9534.700000000001
$ws.Range("K126").Value = 50702.50199999999
$ws.Range("L126").Value = 28604.1
$ws.Range("M126").Value = -48232.50199999999
$ws.Range("N126").Value = -33544.10000000001
$ws.Range("H134").Value = 47803.668
$ws.Range("J134").Value = 47803.668
$ws.Range("L134").Value = 143411.004
$ws.Range("N134").Value = -148481.004
$ws.Range("H136").Value = 128376.5
$ws.Range("J136").Value = 128376.5
$ws.Range("L136").Value = 385129.5
$ws.Range("N136").Value = -390229.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 3988.875
$ws.Range("I16").Value = 1778.1
$ws.Range("J16").Value = 7673.5
$ws.Range("K16").Value = 1778.1
$ws.Range("L16").Value = 7673.5
$ws.Range("M16").Value = -1608.1
$ws.Range("N16").Value = -8013.5
$ws.Range("H40").Value = 22759.045
$ws.Range("I40").Value = 24983.666
$ws.Range("J40").Value = 12748.25
$ws.Range("K40").Value = 24983.666
$ws.Range("L40").Value = 12748.25
$ws.Range("M40").Value = -24847.666
$ws.Range("N40").Value = -13020.25
$ws.Range("H61").Value = 9180.210999999999
$ws.Range("J61").Value = 13700
$ws.Range("L61").Value = 13700
$ws.Range("N61").Value = -14104
$ws.Range("H113").Value = 9180.210999999999
$ws.Range("J113").Value = 13700
$ws.Range("L113").Value = 13700
$ws.Range("N113").Value = -18040

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H39").Value = 20000
$ws.Range("I39").Value = 20000
$ws.Range("K39").Value = 20000
$ws.Range("M39").Value = -19587
$ws.Range("H81").Value = 16221
$ws.Range("I81").Value = 30319.8
$ws.Range("J81").Value = 4472
$ws.Range("K81").Value = 60639.6
$ws.Range("L81").Value = 8944
$ws.Range("M81").Value = -59578.6
$ws.Range("N81").Value = -11066
$ws.Range("H84").Value = 16221
$ws.Range("I84").Value = 30319.8
$ws.Range("J84").Value = 4472
$ws.Range("K84").Value = 303198
$ws.Range("L84").Value = 44720
$ws.Range("M84").Value = -297894
$ws.Range("N84").Value = -55328
$ws.Range("H122").Value = 8382.5
$ws.Range("J122").Value = 10715.308
$ws.Range("L122").Value = 32145.924
$ws.Range("N122").Value = -37045.924
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()
